$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.307.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "'1.865.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'0.7042"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'238.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.07731"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.75%  "
$ws.Range("D9").Value = "'0.3054"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").Value = "'24.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.23%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "'1.874.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "'5.235"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").Value = "'0.7170"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "'89.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'29.404.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "'5.800"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "'240.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "'0.000007803"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "'13.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").Value = "'2.117.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "'7.632"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.30%  "
$ws.Range("D25").Value = "'162.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "'8.935"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "'1.919"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.54%  "
$ws.Range("D30").Value = "'1.369"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.53%  "
$ws.Range("D31").Value = "'1.476"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").Value = "'4.305"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.79%  "
$ws.Range("D33").Value = "'4.040"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").Value = "'0.05209"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "'1.186"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("D36").Value = "'0.7155"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("D37").Value = "'1.003"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "'2.680"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "'2.698"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("D41").Value = "'1.174.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.05%  "
$ws.Range("D42").Value = "'0.9154"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").Value = "'5.986"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "'71.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "'0.4267"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").Value = "'1.002"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "'102.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "'0.5366"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").Value = "'1.749"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").Value = "'9.182"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "'7.002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
